$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (shown in A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 08:41"

# --- Ucrania moves above Israel (row 27/28 swap with refreshed Ucrania data) ---
$ws.Range("A27").Value = "Ucrania"
$ws.Range("B27").Value = 156797
$ws.Range("C27").Value = 2462
$ws.Range("D27").Value = 69543
$ws.Range("E27").Value = 84043
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 33
$ws.Range("H27").Value = 3211

$ws.Range("A28").Value = "Israel"
$ws.Range("B28").Value = 156596
$ws.Range("C28").Value = 992
$ws.Range("D28").Value = 115122
$ws.Range("E28").Value = 40355
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 1119

# --- Standalone daily updates (no reordering) ---
$ws.Range("B60").Value = 47620
$ws.Range("C60").Value = 333
$ws.Range("E60").Value = 3224
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 394

$ws.Range("D75").Value = 18410
$ws.Range("E75").Value = 7730
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 788

# --- Georgia moves above Malta (rows 145-149 shift with refreshed Georgia data) ---
$ws.Range("A145").Value = "Georgia"
$ws.Range("B145").Value = 2392
$ws.Range("C145").Value = 165
$ws.Range("D145").Value = 1369
$ws.Range("E145").Value = 1004
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 19

$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 2352
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 1872
$ws.Range("E146").Value = 465
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 15

$ws.Range("A147").Value = "Guinea-Bisau"
$ws.Range("B147").Value = 2275
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 1127
$ws.Range("E147").Value = 1109
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 39

$ws.Range("A148").Value = "Benin"
$ws.Range("B148").Value = 2267
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 1942
$ws.Range("E148").Value = 285
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 40

$ws.Range("A149").Value = "Botsuana"
$ws.Range("B149").Value = 2252
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 546
$ws.Range("E149").Value = 1696
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 10
